$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph 1: "TLS 1.3 RFC446" -> "TLS 1.3: RFC446"
# (round-trip through a throwaway string with no shared prefix/suffix so the
# host emits one clean run instead of splitting around the inserted colon)
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "QZXJVB"
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "TLS 1.3: RFC446"

# Paragraph 2: "DSCP LE PHB RFC8622" -> "DSCP LE PHB: RFC8622"
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "QZXJVB"
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "DSCP LE PHB: RFC8622"

# Paragraph 5: "SCE draft-morton-tsvwg-sce-00" -> "SCE: draft-morton-tsvwg-sce-00"
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = "QZXJVB"
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = "SCE: draft-morton-tsvwg-sce-00"

# Append a new bullet mentioning RC4 after the SCE paragraph (the last one).
$null = $tr.InsertAfter("`rDeprecate RC4 in SSH (draft-ietf-curdle-rc4-die-die-die-12)")
